# Scheduled Sheets runner: refresh cached Universalis market-price / Leve
# profit figures (currentAveragePrice*, LevePrice*, LeveProfit* columns)
# across all eight Disciple of the Hand tables (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR). Values below are the latest snapshot; cells that did not
# previously exist (e.g. a row's HQ profit column) are created, and columns
# that are no longer applicable for a row are cleared.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1764.2
$ws.Range("J17").Value = 2307.25
$ws.Range("L17").Value = 6921.75
$ws.Range("N17").Value = -7257.75
$ws.Range("H70").Value = 5481002
$ws.Range("I70").Value = 6679222
$ws.Range("J70").Value = 4282781.5
$ws.Range("K70").Value = 20037666
$ws.Range("L70").Value = 12848344.5
$ws.Range("M70").Value = -20037396
$ws.Range("N70").Value = -12848884.5
$ws.Range("H73").Value = 5481002
$ws.Range("I73").Value = 6679222
$ws.Range("J73").Value = 4282781.5
$ws.Range("K73").Value = 20037666
$ws.Range("L73").Value = 12848344.5
$ws.Range("M73").Value = -20036730
$ws.Range("N73").Value = -12850216.5
$ws.Range("H92").Value = 620.4194
$ws.Range("I92").Value = 553.25
$ws.Range("J92").Value = 850.7143
$ws.Range("K92").Value = 553.25
$ws.Range("L92").Value = 850.7143
$ws.Range("M92").Value = 694.75
$ws.Range("N92").Value = -3346.7143
$ws.Range("H98").Value = 1007.36664
$ws.Range("J98").Value = 14995
$ws.Range("L98").Value = 14995
$ws.Range("N98").Value = -17991
$ws.Range("H101").Value = 671.1667
$ws.Range("I101").Value = 609
$ws.Range("J101").Value = 733.3333
$ws.Range("K101").Value = 1827
$ws.Range("L101").Value = 2199.9999
$ws.Range("M101").Value = -205
$ws.Range("N101").Value = -5443.9999
$ws.Range("H106").Value = 3287.5
$ws.Range("I106").Value = 2514.2856
$ws.Range("K106").Value = 2514.2856
$ws.Range("M106").Value = -1883.2856
$ws.Range("H116").Value = 4995.7856
$ws.Range("I116").Value = 3696.7
$ws.Range("J116").Value = 8243.5
$ws.Range("K116").Value = 3696.7
$ws.Range("L116").Value = 8243.5
$ws.Range("M116").Value = -254.6999999999998
$ws.Range("N116").Value = -15127.5
$ws.Range("H122").Value = 1007.36664
$ws.Range("J122").Value = 14995
$ws.Range("L122").Value = 44985
$ws.Range("N122").Value = -49885
$ws.Range("H138").Value = 2261.7837
$ws.Range("I138").Value = 1764.5
$ws.Range("J138").Value = 4393
$ws.Range("K138").Value = 5293.5
$ws.Range("L138").Value = 13179
$ws.Range("M138").Value = -153.5
$ws.Range("N138").Value = -23459

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2686.08
$ws.Range("I32").Value = 2686.08
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 2686.08
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2399.08
$ws.Range("N32").ClearContents()
$ws.Range("H45").Value = 4245.706
$ws.Range("I45").Value = 1399.8334
$ws.Range("J45").Value = 5798
$ws.Range("K45").Value = 1399.8334
$ws.Range("L45").Value = 5798
$ws.Range("M45").Value = -1022.8334
$ws.Range("N45").Value = -6552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 34000
$ws.Range("J55").Value = 34000
$ws.Range("L55").Value = 34000
$ws.Range("N55").Value = -34546
$ws.Range("H99").Value = 3255
$ws.Range("I99").Value = 2004.5
$ws.Range("K99").Value = 2004.5
$ws.Range("M99").Value = -506.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 12690.818
$ws.Range("I62").Value = 3700
$ws.Range("J62").Value = 16062.375
$ws.Range("K62").Value = 3700
$ws.Range("L62").Value = 16062.375
$ws.Range("M62").Value = -3076
$ws.Range("N62").Value = -17310.375
$ws.Range("H65").Value = 12690.818
$ws.Range("I65").Value = 3700
$ws.Range("J65").Value = 16062.375
$ws.Range("K65").Value = 18500
$ws.Range("L65").Value = 80311.875
$ws.Range("M65").Value = -15380
$ws.Range("N65").Value = -86551.875
$ws.Range("H132").Value = 201270.67
$ws.Range("I132").Value = 201270.67
$ws.Range("K132").Value = 603812.01
$ws.Range("M132").Value = -601282.01
$ws.Range("H134").Value = 2203.7827
$ws.Range("I134").Value = 1893.7222
$ws.Range("K134").Value = 5681.1666
$ws.Range("M134").Value = -3146.1666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 944.7
$ws.Range("I114").Value = 638.25
$ws.Range("J114").Value = 1149
$ws.Range("K114").Value = 1914.75
$ws.Range("L114").Value = 3447
$ws.Range("M114").Value = 1339.25
$ws.Range("N114").Value = -9955
$ws.Range("H121").Value = 127800.22
$ws.Range("J121").Value = 8493.25
$ws.Range("L121").Value = 25479.75
$ws.Range("N121").Value = -28099.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 5599
$ws.Range("J23").Value = 5599
$ws.Range("L23").Value = 5599
$ws.Range("N23").Value = -6045
$ws.Range("H44").Value = 5026833
$ws.Range("J44").Value = 5026833
$ws.Range("L44").Value = 5026833
$ws.Range("N44").Value = -5028025
$ws.Range("H102").Value = 17478.94
$ws.Range("J102").Value = 3546.7778
$ws.Range("L102").Value = 3546.7778
$ws.Range("N102").Value = -6790.7778
$ws.Range("H122").Value = 5192
$ws.Range("I122").Value = 5057.3335
$ws.Range("K122").Value = 15172.0005
$ws.Range("M122").Value = -12722.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1113.174
$ws.Range("I22").Value = 715.7692
$ws.Range("J22").Value = 1629.8
$ws.Range("K22").Value = 715.7692
$ws.Range("L22").Value = 1629.8
$ws.Range("M22").Value = -420.7692
$ws.Range("N22").Value = -2219.8
$ws.Range("H27").Value = 1113.174
$ws.Range("I27").Value = 715.7692
$ws.Range("J27").Value = 1629.8
$ws.Range("K27").Value = 715.7692
$ws.Range("L27").Value = 1629.8
$ws.Range("M27").Value = -608.7692
$ws.Range("N27").Value = -1843.8
$ws.Range("H32").Value = 46671.668
$ws.Range("J32").Value = 52507.5
$ws.Range("L32").Value = 52507.5
$ws.Range("N32").Value = -53141.5
$ws.Range("H68").Value = 4676.8887
$ws.Range("I68").Value = 1458.8
$ws.Range("J68").Value = 8699.5
$ws.Range("K68").Value = 1458.8
$ws.Range("L68").Value = 8699.5
$ws.Range("M68").Value = -709.8
$ws.Range("N68").Value = -10197.5
$ws.Range("H71").Value = 4676.8887
$ws.Range("I71").Value = 1458.8
$ws.Range("J71").Value = 8699.5
$ws.Range("K71").Value = 7294
$ws.Range("L71").Value = 43497.5
$ws.Range("M71").Value = -3550
$ws.Range("N71").Value = -50985.5
$ws.Range("H132").Value = 15666.333
$ws.Range("J132").Value = 15000
$ws.Range("L132").Value = 45000
$ws.Range("N132").Value = -50060
$ws.Range("H136").Value = 2532.4666
$ws.Range("I136").Value = 2449.0715
$ws.Range("K136").Value = 7347.2145
$ws.Range("M136").Value = -4797.2145

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6586.5386
$ws.Range("I62").Value = 6036.875
$ws.Range("K62").Value = 6036.875
$ws.Range("M62").Value = -5412.875
$ws.Range("H65").Value = 6586.5386
$ws.Range("I65").Value = 6036.875
$ws.Range("K65").Value = 30184.375
$ws.Range("M65").Value = -27064.375
$ws.Range("H100").Value = 1637.5
$ws.Range("I100").Value = 1814.2858
$ws.Range("K100").Value = 3628.5716
$ws.Range("M100").Value = -3087.5716
$ws.Range("H126").Value = 457071.9
$ws.Range("J126").Value = 1251899.8
$ws.Range("L126").Value = 3755699.4
$ws.Range("N126").Value = -3760639.4
